{"js": "// Collapse the two \"Dataset Non Strutturati\" bullet paragraphs about\n// album / artist info into a single, shortened album bullet:\n//   \"-Informazioni relative ad un album (foto copertina, audio dei brani)\"\n// The old \"-Informazioni relative ad un artista (...)\" paragraph is removed.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst NEW_TEXT = \"-Informazioni relative ad un album (foto copertina, audio dei brani)\";\n\nlet albumParagraph = null;\nlet artistParagraph = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  const t = p.text || \"\";\n  if (t.indexOf(\"Informazioni relative\") === -1) continue;\n  // First bullet: \"-Informazioni relative ad un album (...)\" (originally\n  // split across several runs: \"-Informazioni relative a\" + \"d un\" + \" album ...\").\n  if (t.indexOf(\"un album\") !== -1) {\n    albumParagraph = p;\n  }\n  // Second bullet: \"-Informazioni relative ad un artista (...)\".\n  if (t.indexOf(\"un artista\") !== -1) {\n    artistParagraph = p;\n  }\n}\n\nif (albumParagraph) {\n  // Replace all of this paragraph's runs with the new, shortened text.\n  albumParagraph.insertText(NEW_TEXT, \"Replace\");\n}\n\nif (artistParagraph) {\n  // The artist paragraph is merged away entirely.\n  artistParagraph.delete();\n}\n\nawait context.sync();\n", "ps1": "# Collapse the two \"Dataset Non Strutturati\" bullet paragraphs about\n# album / artist info into a single, shortened album bullet:\n#   \"-Informazioni relative ad un album (foto copertina, audio dei brani)\"\n# The old \"-Informazioni relative ad un artista (...)\" paragraph is removed.\n\n$d = $word.ActiveDocument\n\n$newText = \"-Informazioni relative ad un album (foto copertina, audio dei brani)\"\n\n# First pass: rewrite the album paragraph's text (collapsing its many runs\n# into a single run with the shortened wording).\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*Informazioni relative*un album*\") {\n        $r = $p.Range\n        # Drop the trailing paragraph mark from the range so the assignment\n        # only replaces the paragraph's text, not its break.\n        $r.End = $r.End - 1\n        $r.Text = $newText\n        break\n    }\n}\n\n# Second pass: delete the artist paragraph entirely (re-enumerate since the\n# previous edit may have shifted ranges).\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*Informazioni relative*un artista*\") {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
